# The author regenerated the vocabulary sheet from the Google Sheet source,
# which re-ordered/added a couple of metadata columns in the "term" table.
# Net effect on the OOXML: two new (blank) columns were inserted before the
# old column C, shifting every existing column C..AN two places right to
# E..AP (and the table's dimension grows from A1:AN30 to A1:AP30).
#
# On top of that generic shift, row 23 (the field-header row) had its
# "skos:altLabel" header moved back next to skos:prefLabel and a brand new
# "skos:notation" header introduced, row 24 gained a "task" subject-type
# value, row 25 gained a "var" value, and the auto-updated "dct:modified"
# timestamp in B21 ticked forward.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two blank columns before column C - this reproduces the bulk of
#    the diff: every populated cell from old-C onward slides two columns to
#    the right (C->E, D->F, E->G, ... AN->AP), and the sheet's used range
#    grows to include the two new trailing blank columns (AO:AP).
$ws.Range("C:D").Insert()

# 2) Row 23 header row: after the generic shift, the old G23 ("skos:altLabel...")
#    landed in I23. The regenerated sheet instead keeps it right after column B,
#    and introduces a brand-new "skos:notation" header in the column next to it.
$ws.Range("C23").Value = 'skos:altLabel(separator=",")'
$ws.Range("D23").Value = "skos:notation"
$ws.Range("I23").Value = ""

# 3) Row 24 ("subject" term row) gained a new value in column C.
$ws.Range("C24").Value = "task"

# 4) Row 25 ("variable" term row) gained a new value in column D.
$ws.Range("D25").Value = "var"

# 5) The auto-generated "dct:modified" timestamp advanced.
$ws.Range("B21").Value = "2023-09-13T14:52:23+00:00"
